# The commit swaps the presentation's theme palette from the "Integral"
# theme to the default Office theme (ppt/theme/theme1.xml, the theme
# actually applied to the slide master / slides) while the companion
# theme used only by the notes master (ppt/theme/theme2.xml) keeps the
# generic "Office Theme" look it always had.
#
# PowerPoint's COM ColorFormat.RGB uses the classic OLE_COLOR packing
# (0x00BBGGRR), i.e. the byte order is reversed relative to the "RRGGBB"
# hex strings you see in the theme XML, so build each value from its
# R/G/B components instead of hard-coding the swapped decimal number.
function ToOleColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office theme "Office" color scheme (the target palette for theme1.xml).
$officeColors = @(
    @(0x00, 0x00, 0x00), # 1  dk1
    @(0xFF, 0xFF, 0xFF), # 2  lt1
    @(0x44, 0x54, 0x6A), # 3  dk2
    @(0xE7, 0xE6, 0xE6), # 4  lt2
    @(0x5B, 0x9B, 0xD5), # 5  accent1
    @(0xED, 0x7D, 0x31), # 6  accent2
    @(0xA5, 0xA5, 0xA5), # 7  accent3
    @(0xFF, 0xC0, 0x00), # 8  accent4
    @(0x44, 0x72, 0xC4), # 9  accent5
    @(0x70, 0xAD, 0x47), # 10 accent6
    @(0x05, 0x63, 0xC1), # 11 hlink
    @(0x95, 0x4F, 0x72)  # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $rgb = $officeColors[$i - 1]
    $colorScheme.Colors($i).RGB = ToOleColor $rgb[0] $rgb[1] $rgb[2]
}
